$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 1.277433666666667
$ws.Cells.Item(2, 8).Value = 3.832301
$ws.Cells.Item(2, 9).Value = 0.01786062203930835
$ws.Cells.Item(2, 10).Value = 0.01786062203930835
$ws.Cells.Item(2, 13).Value = 42.09975866666667
$ws.Cells.Item(2, 14).Value = 126.299276
$ws.Cells.Item(2, 15).Value = 0.3315552933456474
$ws.Cells.Item(2, 16).Value = 0.3315552933456474
$ws.Cells.Item(2, 17).Value = 53.77964907934177
$ws.Cells.Item(2, 18).Value = 484.016841714076
$ws.Cells.Item(2, 19).Value = 0.005921783779578616
$ws.Cells.Item(2, 20).Value = 0.005921783779578615

$ws.Cells.Item(3, 7).Value = 1.277433666666667
$ws.Cells.Item(3, 8).Value = 3.832301
$ws.Cells.Item(3, 9).Value = 0.01786062203930835
$ws.Cells.Item(3, 10).Value = 0.01786062203930835
$ws.Cells.Item(3, 15).Value = 0.4502223747274475
$ws.Cells.Item(3, 16).Value = 0.4502223747274475
$ws.Cells.Item(3, 17).Value = 73.02794377427756
$ws.Cells.Item(3, 18).Value = 657.251493968498
$ws.Cells.Item(3, 19).Value = 0.008041251668646793
$ws.Cells.Item(3, 20).Value = 0.008041251668646793

$ws.Cells.Item(4, 7).Value = 1.277433666666667
$ws.Cells.Item(4, 8).Value = 3.832301
$ws.Cells.Item(4, 9).Value = 0.01786062203930835
$ws.Cells.Item(4, 10).Value = 0.01786062203930835
$ws.Cells.Item(4, 13).Value = 27.596267
$ws.Cells.Item(4, 14).Value = 82.78880100000001
$ws.Cells.Item(4, 15).Value = 0.2173335118824389
$ws.Cells.Item(4, 16).Value = 0.2173335118824389
$ws.Cells.Item(4, 17).Value = 35.25240054012233
$ws.Cells.Item(4, 18).Value = 317.271604861101
$ws.Cells.Item(4, 19).Value = 0.003881711712207771
$ws.Cells.Item(4, 20).Value = 0.003881711712207771

$ws.Cells.Item(5, 7).Value = 1.277433666666667
$ws.Cells.Item(5, 8).Value = 3.832301
$ws.Cells.Item(5, 9).Value = 0.01786062203930835
$ws.Cells.Item(5, 10).Value = 0.01786062203930835
$ws.Cells.Item(5, 13).Value = 0.1128593333333333
$ws.Cells.Item(5, 14).Value = 0.338578
$ws.Cells.Item(5, 15).Value = 0.0008888200444663087
$ws.Cells.Item(5, 16).Value = 0.0008888200444663087
$ws.Cells.Item(5, 17).Value = 0.1441703119975555
$ws.Cells.Item(5, 18).Value = 1.297532807978
$ws.Cells.Item(5, 19).Value = 0.00001587487887517399
$ws.Cells.Item(5, 20).Value = 0.00001587487887517398

$ws.Cells.Item(6, 7).Value = 45.44725166666667
$ws.Cells.Item(6, 9).Value = 0.6354272679079697
$ws.Cells.Item(6, 10).Value = 0.6354272679079697
$ws.Cells.Item(6, 13).Value = 42.09975866666667
$ws.Cells.Item(6, 14).Value = 126.299276
$ws.Cells.Item(6, 15).Value = 0.3315552933456474
$ws.Cells.Item(6, 16).Value = 0.3315552933456474
$ws.Cells.Item(6, 17).Value = 1913.318327229931
$ws.Cells.Item(6, 18).Value = 17219.86494506938
$ws.Cells.Item(6, 19).Value = 0.2106792742110501
$ws.Cells.Item(6, 20).Value = 0.2106792742110501

$ws.Cells.Item(7, 7).Value = 45.44725166666667
$ws.Cells.Item(7, 9).Value = 0.6354272679079697
$ws.Cells.Item(7, 10).Value = 0.6354272679079697
$ws.Cells.Item(7, 15).Value = 0.4502223747274475
$ws.Cells.Item(7, 16).Value = 0.4502223747274475
$ws.Cells.Item(7, 17).Value = 2598.114818806333
$ws.Cells.Item(7, 19).Value = 0.2860835735241001
$ws.Cells.Item(7, 20).Value = 0.2860835735241001

$ws.Cells.Item(8, 7).Value = 45.44725166666667
$ws.Cells.Item(8, 9).Value = 0.6354272679079697
$ws.Cells.Item(8, 10).Value = 0.6354272679079697
$ws.Cells.Item(8, 13).Value = 27.596267
$ws.Cells.Item(8, 14).Value = 82.78880100000001
$ws.Cells.Item(8, 15).Value = 0.2173335118824389
$ws.Cells.Item(8, 16).Value = 0.2173335118824389
$ws.Cells.Item(8, 17).Value = 1254.174491409528
$ws.Cells.Item(8, 18).Value = 11287.57042268576
$ws.Cells.Item(8, 19).Value = 0.1380996396803024
$ws.Cells.Item(8, 20).Value = 0.1380996396803024

$ws.Cells.Item(9, 7).Value = 45.44725166666667
$ws.Cells.Item(9, 9).Value = 0.6354272679079697
$ws.Cells.Item(9, 10).Value = 0.6354272679079697
$ws.Cells.Item(9, 13).Value = 0.1128593333333333
$ws.Cells.Item(9, 14).Value = 0.338578
$ws.Cells.Item(9, 15).Value = 0.0008888200444663087
$ws.Cells.Item(9, 16).Value = 0.0008888200444663087
$ws.Cells.Item(9, 17).Value = 5.129146524932222
$ws.Cells.Item(9, 18).Value = 46.16231872439
$ws.Cells.Item(9, 19).Value = 0.0005647804925170666
$ws.Cells.Item(9, 20).Value = 0.0005647804925170666

$ws.Cells.Item(10, 7).Value = 23.96074166666667
$ws.Cells.Item(10, 8).Value = 71.88222500000001
$ws.Cells.Item(10, 9).Value = 0.3350105464235513
$ws.Cells.Item(10, 10).Value = 0.3350105464235513
$ws.Cells.Item(10, 13).Value = 42.09975866666667
$ws.Cells.Item(10, 14).Value = 126.299276
$ws.Cells.Item(10, 15).Value = 0.3315552933456474
$ws.Cells.Item(10, 16).Value = 0.3315552933456474
$ws.Cells.Item(10, 17).Value = 1008.741441641011
$ws.Cells.Item(10, 18).Value = 9078.672974769101
$ws.Cells.Item(10, 19).Value = 0.1110745199933462
$ws.Cells.Item(10, 20).Value = 0.1110745199933462

$ws.Cells.Item(11, 7).Value = 23.96074166666667
$ws.Cells.Item(11, 8).Value = 71.88222500000001
$ws.Cells.Item(11, 9).Value = 0.3350105464235513
$ws.Cells.Item(11, 10).Value = 0.3350105464235513
$ws.Cells.Item(11, 15).Value = 0.4502223747274475
$ws.Cells.Item(11, 16).Value = 0.4502223747274475
$ws.Cells.Item(11, 17).Value = 1369.780475403672
$ws.Cells.Item(11, 18).Value = 12328.02427863305
$ws.Cells.Item(11, 19).Value = 0.1508292437695511
$ws.Cells.Item(11, 20).Value = 0.1508292437695511

$ws.Cells.Item(12, 7).Value = 23.96074166666667
$ws.Cells.Item(12, 8).Value = 71.88222500000001
$ws.Cells.Item(12, 9).Value = 0.3350105464235513
$ws.Cells.Item(12, 10).Value = 0.3350105464235513
$ws.Cells.Item(12, 13).Value = 27.596267
$ws.Cells.Item(12, 14).Value = 82.78880100000001
$ws.Cells.Item(12, 15).Value = 0.2173335118824389
$ws.Cells.Item(12, 16).Value = 0.2173335118824389
$ws.Cells.Item(12, 17).Value = 661.2270245513583
$ws.Cells.Item(12, 18).Value = 5951.043220962226
$ws.Cells.Item(12, 19).Value = 0.07280901857188522
$ws.Cells.Item(12, 20).Value = 0.07280901857188522

$ws.Cells.Item(13, 7).Value = 23.96074166666667
$ws.Cells.Item(13, 8).Value = 71.88222500000001
$ws.Cells.Item(13, 9).Value = 0.3350105464235513
$ws.Cells.Item(13, 10).Value = 0.3350105464235513
$ws.Cells.Item(13, 13).Value = 0.1128593333333333
$ws.Cells.Item(13, 14).Value = 0.338578
$ws.Cells.Item(13, 15).Value = 0.0008888200444663087
$ws.Cells.Item(13, 16).Value = 0.0008888200444663087
$ws.Cells.Item(13, 17).Value = 2.704193330672222
$ws.Cells.Item(13, 18).Value = 24.33773997605
$ws.Cells.Item(13, 19).Value = 0.0002977640887688632
$ws.Cells.Item(13, 20).Value = 0.0002977640887688632

$ws.Cells.Item(14, 7).Value = 0.8369233333333334
$ws.Cells.Item(14, 8).Value = 2.51077
$ws.Cells.Item(14, 9).Value = 0.01170156362917063
$ws.Cells.Item(14, 10).Value = 0.01170156362917063
$ws.Cells.Item(14, 13).Value = 42.09975866666667
$ws.Cells.Item(14, 14).Value = 126.299276
$ws.Cells.Item(14, 15).Value = 0.3315552933456474
$ws.Cells.Item(14, 16).Value = 0.3315552933456474
$ws.Cells.Item(14, 17).Value = 35.23427035583556
$ws.Cells.Item(14, 18).Value = 317.10843320252
$ws.Cells.Item(14, 19).Value = 0.003879715361672427
$ws.Cells.Item(14, 20).Value = 0.003879715361672427

$ws.Cells.Item(15, 7).Value = 0.8369233333333334
$ws.Cells.Item(15, 8).Value = 2.51077
$ws.Cells.Item(15, 9).Value = 0.01170156362917063
$ws.Cells.Item(15, 10).Value = 0.01170156362917063
$ws.Cells.Item(15, 15).Value = 0.4502223747274475
$ws.Cells.Item(15, 16).Value = 0.4502223747274475
$ws.Cells.Item(15, 17).Value = 47.84498148505111
$ws.Cells.Item(15, 18).Value = 430.60483336546
$ws.Cells.Item(15, 19).Value = 0.005268305765149531
$ws.Cells.Item(15, 20).Value = 0.005268305765149531

$ws.Cells.Item(16, 7).Value = 0.8369233333333334
$ws.Cells.Item(16, 8).Value = 2.51077
$ws.Cells.Item(16, 9).Value = 0.01170156362917063
$ws.Cells.Item(16, 10).Value = 0.01170156362917063
$ws.Cells.Item(16, 13).Value = 27.596267
$ws.Cells.Item(16, 14).Value = 82.78880100000001
$ws.Cells.Item(16, 15).Value = 0.2173335118824389
$ws.Cells.Item(16, 16).Value = 0.2173335118824389
$ws.Cells.Item(16, 17).Value = 23.09595976519667
$ws.Cells.Item(16, 18).Value = 207.86363788677
$ws.Cells.Item(16, 19).Value = 0.00254314191804347
$ws.Cells.Item(16, 20).Value = 0.00254314191804347

$ws.Cells.Item(17, 7).Value = 0.8369233333333334
$ws.Cells.Item(17, 8).Value = 2.51077
$ws.Cells.Item(17, 9).Value = 0.01170156362917063
$ws.Cells.Item(17, 10).Value = 0.01170156362917063
$ws.Cells.Item(17, 13).Value = 0.1128593333333333
$ws.Cells.Item(17, 14).Value = 0.338578
$ws.Cells.Item(17, 15).Value = 0.0008888200444663087
$ws.Cells.Item(17, 16).Value = 0.0008888200444663087
$ws.Cells.Item(17, 17).Value = 0.0944546094511111
$ws.Cells.Item(17, 18).Value = 0.85009148506
$ws.Cells.Item(17, 19).Value = 0.00001040058430520478
$ws.Cells.Item(17, 20).Value = 0.00001040058430520478
